$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "/home/daniel/Spike Data/Matlab files/Exp 18 baseline data.mat"
$ws.Range("B7").Value = 2
$ws.Range("C7").Value = 3.97
$ws.Range("D7").Value = 3
$ws.Range("E7").Value = 15440
$ws.Range("F7").Value = 19220
$ws.Range("G7").Value = 62700
$ws.Range("H7").Value = 72460
$ws.Range("I7").Value = 156300
$ws.Range("J7").Value = 167500

$ws.Range("A8").Value = "/home/daniel/Spike Data/Matlab files/exp 28 baseline.mat"
$ws.Range("B8").Value = 2
$ws.Range("C8").Value = 2.77
$ws.Range("D8").Value = 3
$ws.Range("E8").Value = 7537
$ws.Range("F8").Value = 15450
$ws.Range("G8").Value = 270600
$ws.Range("H8").Value = 275400
$ws.Range("I8").Value = 434000
$ws.Range("J8").Value = 440000

$ws.Range("J16").Select()
